$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original row 1 values (as text) before overwriting them.
$oldA1 = $ws.Range("A1").Text
$oldB1 = $ws.Range("B1").Text

# Row 1 now holds the new user/password pair.
$ws.Range("A1").Value2 = "Nick"
$ws.Range("B1").Value2 = "Password"

# The previous row 1 contents ("Tomek" / "123") are appended as a new row 4,
# keeping the numeric-looking "123" stored as text.
$ws.Range("A4").Value2 = $oldA1
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value2 = $oldB1
